$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.767.31'
$ws.Range("E2").Value = '  -1.19%  '

$ws.Range("D3").Value = '2.445.39'
$ws.Range("E3").Value = '  -0.10%  '

$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.95'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.81%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.86'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.99%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.532'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.60%  '

$ws.Range("D9").Value = '2.436.69'
$ws.Range("E9").Value = '  -0.38%  '

$ws.Range("E10").Value = '  +2.07%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.16'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.03%  '

$ws.Range("E13").Value = '  -1.83%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.87'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.35%  '

$ws.Range("D15").Value = '2.888.56'
$ws.Range("E15").Value = '  +1.00%  '

$ws.Range("E16").Value = '  -1.16%  '

$ws.Range("D17").Value = '61.703.99'
$ws.Range("E17").Value = '  -0.76%  '

$ws.Range("D18").Value = '2.451.09'
$ws.Range("E18").Value = '  +0.51%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.60'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.81%  '

$ws.Range("E20").Value = '  +1.33%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '325.04'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.21%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.07'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.24%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.06'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.13%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.94'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.32%  '

$ws.Range("E25").Value = '  +0.10%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '64.89'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.52%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.12'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.22%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '581.17'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -9.17%  '

$ws.Range("D29").Value = '2.567.94'
$ws.Range("E29").Value = '  +0.52%  '

$ws.Range("E30").Value = '  +0.10%  '

$ws.Range("E31").Value = '  -3.50%  '

$ws.Range("E32").Value = '  -1.98%  '

$ws.Range("E33").Value = '  -5.44%  '

$ws.Range("E34").Value = '  -1.32%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.132'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.98%  '

$ws.Range("E36").Value = '  +0.15%  '

$ws.Range("E37").Value = '  -5.58%  '

$ws.Range("E38").Value = '  -1.14%  '

$ws.Range("E39").Value = '  -3.96%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '151.17'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.28%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '18.26'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.92%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.13'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.63%  '

$ws.Range("E43").Value = '  +0.00%  '

$ws.Range("E44").Value = '  -4.06%  '

$ws.Range("E45").Value = '  -2.63%  '

$ws.Range("E46").Value = '  -5.78%  '

$ws.Range("D47").Value = '0.0₆0285'
$ws.Range("E47").Value = '  +22.34%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '142.57'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.80%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.56'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.88%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.599'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.15%  '

$ws.Range("E51").Value = '  -0.35%  '
